$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Agregar cta 60109 al grupo de clientes en B5
$ws.Range("B5").Value2 = $ws.Range("B5").Value2 + ".60109"

# Actualizar la vista: seleccionar B6 y volver a A1 como celda superior izquierda
$ws.Range("B6").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
